$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 100
$ws.Range("I5").Value = 100
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 100
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 15
$ws.Range("N5").ClearContents()
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").ClearContents()
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("H52").Value = 400100
$ws.Range("J52").Value = 500000
$ws.Range("L52").Value = 1500000
$ws.Range("N52").Value = -1500320
$ws.Range("H111").Value = 5170.0835
$ws.Range("I111").Value = 5286.3335
$ws.Range("J111").Value = 4821.3335
$ws.Range("K111").Value = 15859.0005
$ws.Range("L111").Value = 14464.0005
$ws.Range("M111").Value = -12792.0005
$ws.Range("N111").Value = -20598.0005
$ws.Range("H112").Value = 1387.3143
$ws.Range("I112").Value = 353
$ws.Range("J112").Value = 1559.7
$ws.Range("K112").Value = 1059
$ws.Range("L112").Value = 4679.1
$ws.Range("M112").Value = 49
$ws.Range("N112").Value = -6895.1
$ws.Range("H113").Value = 2736.4
$ws.Range("I113").Value = 2000
$ws.Range("J113").Value = 3052
$ws.Range("K113").Value = 2000
$ws.Range("L113").Value = 3052
$ws.Range("M113").Value = 1254
$ws.Range("N113").Value = -9560
$ws.Range("H137").Value = 4381.9
$ws.Range("I137").Value = 4340
$ws.Range("J137").Value = 4390.28
$ws.Range("K137").Value = 13020
$ws.Range("L137").Value = 13170.84
$ws.Range("M137").Value = -10470
$ws.Range("N137").Value = -18270.84

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8665.088
$ws.Range("I32").Value = 7542.1973
$ws.Range("J32").Value = 30000
$ws.Range("K32").Value = 7542.1973
$ws.Range("L32").Value = 30000
$ws.Range("M32").Value = -7255.1973
$ws.Range("N32").Value = -30574
$ws.Range("H88").Value = 5890.769
$ws.Range("I88").Value = 18002
$ws.Range("J88").Value = 2257.4
$ws.Range("K88").Value = 18002
$ws.Range("L88").Value = 2257.4
$ws.Range("M88").Value = -17596
$ws.Range("N88").Value = -3069.4
$ws.Range("H91").Value = 5890.769
$ws.Range("I91").Value = 18002
$ws.Range("J91").Value = 2257.4
$ws.Range("K91").Value = 18002
$ws.Range("L91").Value = 2257.4
$ws.Range("M91").Value = -16598
$ws.Range("N91").Value = -5065.4
$ws.Range("H131").Value = 50000
$ws.Range("J131").Value = 50000
$ws.Range("L131").Value = 50000
$ws.Range("N131").Value = -60080

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H125").Value = 79645
$ws.Range("J125").Value = 79645
$ws.Range("L125").Value = 79645
$ws.Range("N125").Value = -89485

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 100
$ws.Range("I6").Value = 100
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 100
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = 13
$ws.Range("N6").ClearContents()
$ws.Range("H16").Value = 1421.3846
$ws.Range("I16").Value = 1466.1428
$ws.Range("J16").Value = 1369.1666
$ws.Range("K16").Value = 1466.1428
$ws.Range("L16").Value = 1369.1666
$ws.Range("M16").Value = -1179.1428
$ws.Range("N16").Value = -1943.1666
$ws.Range("H22").Value = 251.61539
$ws.Range("I22").Value = 207.1
$ws.Range("K22").Value = 207.1
$ws.Range("M22").Value = 142.9
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("M25").ClearContents()
$ws.Range("N25").ClearContents()
$ws.Range("H31").Value = 8597.081
$ws.Range("I31").Value = 3249.7896
$ws.Range("J31").Value = 14241.444
$ws.Range("K31").Value = 3249.7896
$ws.Range("L31").Value = 14241.444
$ws.Range("M31").Value = -2954.7896
$ws.Range("N31").Value = -14831.444
$ws.Range("H34").Value = 8597.081
$ws.Range("I34").Value = 3249.7896
$ws.Range("J34").Value = 14241.444
$ws.Range("K34").Value = 3249.7896
$ws.Range("L34").Value = 14241.444
$ws.Range("M34").Value = -3047.7896
$ws.Range("N34").Value = -14645.444
$ws.Range("H50").Value = 20092
$ws.Range("J50").Value = 20092
$ws.Range("L50").Value = 20092
$ws.Range("N50").Value = -21342
$ws.Range("H51").Value = 22253.818
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 22253.818
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 22253.818
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = -23725.818
$ws.Range("H60").Value = 13199.8
$ws.Range("I60").Value = 13199.8
$ws.Range("J60").Value = 0
$ws.Range("K60").Value = 13199.8
$ws.Range("L60").Value = 0
$ws.Range("M60").Value = -12688.8
$ws.Range("N60").ClearContents()
$ws.Range("H61").Value = 22253.818
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 22253.818
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 22253.818
$ws.Range("M61").ClearContents()
$ws.Range("N61").Value = -22949.818
$ws.Range("H113").Value = 1421.3846
$ws.Range("I113").Value = 1466.1428
$ws.Range("J113").Value = 1369.1666
$ws.Range("K113").Value = 1466.1428
$ws.Range("L113").Value = 1369.1666
$ws.Range("M113").Value = 703.8571999999999
$ws.Range("N113").Value = -5709.1666

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 872.6875
$ws.Range("I23").Value = 2700.75
$ws.Range("J23").Value = 263.33334
$ws.Range("K23").Value = 8102.25
$ws.Range("L23").Value = 790.0000200000001
$ws.Range("M23").Value = -7867.25
$ws.Range("N23").Value = -1260.00002
$ws.Range("H34").Value = 3849.2593
$ws.Range("I34").Value = 210
$ws.Range("J34").Value = 4304.1665
$ws.Range("K34").Value = 630
$ws.Range("L34").Value = 12912.4995
$ws.Range("M34").Value = -546
$ws.Range("N34").Value = -13080.4995
$ws.Range("H39").Value = 8873.076999999999
$ws.Range("J39").Value = 8873.076999999999
$ws.Range("L39").Value = 26619.231
$ws.Range("N39").Value = -27207.231
$ws.Range("H51").Value = 2000
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 2000
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 6000
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = -6920
$ws.Range("H55").Value = 8747.5
$ws.Range("I55").Value = 5000
$ws.Range("J55").Value = 9996.666999999999
$ws.Range("K55").Value = 15000
$ws.Range("L55").Value = 29990.001
$ws.Range("M55").Value = -14823
$ws.Range("N55").Value = -30344.001
$ws.Range("H58").Value = 3146.1538
$ws.Range("I58").Value = 2050
$ws.Range("J58").Value = 3345.4546
$ws.Range("K58").Value = 6150
$ws.Range("L58").Value = 10036.3638
$ws.Range("M58").Value = -6022
$ws.Range("N58").Value = -10292.3638
$ws.Range("H107").Value = 2137616
$ws.Range("J107").Value = 1041.9744
$ws.Range("L107").Value = 3125.9232
$ws.Range("N107").Value = -6965.9232
$ws.Range("H110").Value = 0
$ws.Range("I110").Value = 0
$ws.Range("K110").Value = 0
$ws.Range("M110").ClearContents()
$ws.Range("H115").Value = 2212
$ws.Range("J115").Value = 2249.2222
$ws.Range("L115").Value = 6747.6666
$ws.Range("N115").Value = -9097.6666
$ws.Range("H122").Value = 1183
$ws.Range("I122").Value = 630.5714
$ws.Range("J122").Value = 1569.7
$ws.Range("K122").Value = 5675.1426
$ws.Range("L122").Value = 14127.3
$ws.Range("M122").Value = -3225.1426
$ws.Range("N122").Value = -19027.3
$ws.Range("H138").Value = 7532.524
$ws.Range("J138").Value = 3724.923
$ws.Range("L138").Value = 11174.769
$ws.Range("N138").Value = -21454.769

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 44332.406
$ws.Range("I132").Value = 170079.83
$ws.Range("J132").Value = 8404.571
$ws.Range("K132").Value = 510239.49
$ws.Range("L132").Value = 25213.713
$ws.Range("M132").Value = -507709.49
$ws.Range("N132").Value = -30273.713

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("N12").ClearContents()
$ws.Range("H136").Value = 5527.4414
$ws.Range("I136").Value = 3774.923
$ws.Range("J136").Value = 11223.125
$ws.Range("K136").Value = 11324.769
$ws.Range("L136").Value = 33669.375
$ws.Range("M136").Value = -8774.769
$ws.Range("N136").Value = -38769.375

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1585.3125
$ws.Range("I107").Value = 874
$ws.Range("J107").Value = 4667.6665
$ws.Range("K107").Value = 2622
$ws.Range("L107").Value = 14002.9995
$ws.Range("M107").Value = -702
$ws.Range("N107").Value = -17842.9995
$ws.Range("H132").Value = 4327.6665
$ws.Range("I132").Value = 3009
$ws.Range("J132").Value = 6173.8
$ws.Range("K132").Value = 9027
$ws.Range("L132").Value = 18521.4
$ws.Range("M132").Value = -6497
$ws.Range("N132").Value = -23581.4
